$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above the current row 34, shifting existing rows
# 34-100 down to 35-101 (dimension grows from A1:T100 to A1:T101).
$ws.Rows("34:34").Insert()

# Populate the newly inserted row 34 with the new price-report entry.
$ws.Range("A34").Value = 5
$ws.Range("B34").Value = "Macroferia Regional de Talca"
$ws.Range("C34").Value = "Maule"
$ws.Range("D34").Value = 44930
$ws.Range("E34").Value = 7
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100101
$ws.Range("H34").Value = "Berries"
$ws.Range("I34").Value = 100101001
$ws.Range("J34").Value = "Arándano (blue)"
$ws.Range("K34").Value = "Sin especificar"
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 170
$ws.Range("N34").Value = 2800
$ws.Range("O34").Value = 3000
$ws.Range("P34").Value = 2859
$ws.Range("Q34").Value = "$/bandeja 2 kilos"
$ws.Range("R34").Value = "Provincia de Curicó"
$ws.Range("S34").Value = 1430
$ws.Range("T34").Value = 2
